# ------------------------------------------------------------------
# Commit: [ADDITIONAL SCRAPING] added code to scrape more data about a
# player's batting performance in a match, also updated the excel sheets
#
# This script:
#   1. Inserts a new "Player Info" sheet at the front of the workbook
#      with basic player metadata.
#   2. Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and
#      "ODI Bowling", rewriting the full howstat URL values down to the
#      bare numeric match code.
#   3. Appends a new "ODI Batting Extra" sheet with extended per-innings
#      batting stats keyed by MATCH_CODE.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Get-MatchCode($url) {
    $marker = "MatchCode="
    return $url.Substring($url.IndexOf($marker) + $marker.Length)
}

# ------------------------------------------------------------------
# 1. "Player Info" sheet — inserted before the existing sheets.
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $piHeaders.Length; $col++) {
    $playerInfo.Cells.Item(1, $col).Value = $piHeaders[$col - 1]
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4577"
$playerInfo.Cells.Item(2, 2).Value = "Christopher B Sole"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium"

# ------------------------------------------------------------------
# 2. "ODI Batting" — MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$battingLinkCol = 4

$batting.Cells.Item(1, $battingLinkCol).Value = "MATCH_CODE"

$battingLastRow = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, $battingLinkCol)
    $url = $cell.Text
    if ($url) {
        $code = Get-MatchCode $url
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# INNING_NUMBER (col B) was only ever populated for innings the player
# actually batted in; the blank placeholder cells left behind by the
# scraper are dropped entirely now.
$battingEmptyInningRows = @(2, 3, 4, 7, 10, 12, 13, 14, 19, 20, 22)
foreach ($r in $battingEmptyInningRows) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# ------------------------------------------------------------------
# 3. "ODI Bowling" — MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlingLinkCol = 2

$bowling.Cells.Item(1, $bowlingLinkCol).Value = "MATCH_CODE"

$bowlingLastRow = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, $bowlingLinkCol)
    $url = $cell.Text
    if ($url) {
        $code = Get-MatchCode $url
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ------------------------------------------------------------------
# 4. "ODI Batting Extra" — appended after "ODI Bowling".
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $extraHeaders[$col - 1]
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3934", 11,   $null, $null, $null,    "NO"),
    @("3979", 10,   "0",   "0",   "0.98%",  "NO"),
    @("3980", 10,   "0",   "0",   "0.58%",  "NO"),
    @("4048", 11,   $null, $null, $null,    "NO"),
    @("4049", $null,$null, $null, $null,    "NO"),
    @("4078", 10,   "0",   "0",   "1.97%",  "NO"),
    @("4155", 10,   $null, $null, $null,    "NO"),
    @("4158", 10,   "1",   "0",   "2.85%",  "NO"),
    @("4512", 9,    $null, $null, $null,    "NO"),
    @("4576", 10,   $null, $null, $null,    "NO"),
    @("4578", 9,    $null, $null, $null,    "NO"),
    @("4581", 11,   "0",   "0",   "1.75%",  "NO"),
    @("4625", 10,   "0",   "0",   "0.98%",  "NO"),
    @("4629", 10,   "1",   "0",   "1.53%",  "NO"),
    @("4632", 10,   "0",   "0",   "0.39%",  "NO"),
    @("4677", 10,   $null, $null, $null,    "NO"),
    @("4681", 10,   $null, $null, $null,    "YES"),
    @("4680", 10,   "2",   "0",   "7.21%",  "NO"),
    @("4684", 10,   "1",   "0",   "1.82%",  "NO"),
    @("4703", $null,$null, $null, $null,    $null)
)

$rowIdx = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($rowIdx, 1).NumberFormat = "@"
    $extra.Cells.Item($rowIdx, 1).Value = $row[0]

    if ($null -eq $row[1]) {
        $extra.Cells.Item($rowIdx, 2).Value = ""
    } else {
        $extra.Cells.Item($rowIdx, 2).Value = $row[1]
    }

    for ($col = 3; $col -le 5; $col++) {
        $v = $row[$col - 1]
        if ($null -eq $v) {
            $extra.Cells.Item($rowIdx, $col).Value = ""
        } else {
            $extra.Cells.Item($rowIdx, $col).NumberFormat = "@"
            $extra.Cells.Item($rowIdx, $col).Value = $v
        }
    }

    $manOfMatch = $row[5]
    if ($null -eq $manOfMatch) {
        $extra.Cells.Item($rowIdx, 6).Value = ""
    } else {
        $extra.Cells.Item($rowIdx, 6).Value = $manOfMatch
    }

    $rowIdx++
}

# ------------------------------------------------------------------
# Leave the workbook's active-tab selection on the first sheet, same
# as before the edit (activeTab="0").
# ------------------------------------------------------------------
$playerInfo.Activate()
